# "Changes in two sum 2"
# The "Two Sum II - Input Array Is Sorted" row (row 10) had its
# "efficient O" column corrected from O(log(n)) to O(n), and the
# now-unused "O(log(n))" shared string is dropped automatically on save.
# The active selection is also moved onto the edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Value = "O(n)"

[void]$ws.Range("F10").Select()
